$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.27"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'1.08%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'27.19"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'1.12%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'4.707"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'0.95%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.06194"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'3.29%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'6.728"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'0.91%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.8512"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.85%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9143"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-0.84%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1414"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'1.47%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.04560"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'0.83%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07086"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'0.82%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.03132"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-0.03%"
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'-0.92%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001529"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'0.24%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.0006161"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'1.62%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.006049"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.33%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'3.459"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.02%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'3.164"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'0.10%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'2.194"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'0.39%"
$ws.Range("E19").ClearFormats()
$ws.Range("E21").Value = "'1.73%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'4.101"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-1.25%"
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'-0.12%"
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'-0.01%"
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'-5.70%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'0.05%"
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'-6.63%"
$ws.Range("E27").ClearFormats()
$ws.Range("D40").Value = "'0.03942"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'2.63%"
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'-0.18%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'6.99%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.002161"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-10.70%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.01381"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-9.53%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005152"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'0.61%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'0.05%"
$ws.Range("E46").ClearFormats()
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").ClearFormats()
